$wb = $excel.ActiveWorkbook

# --- Fix the "cep290_unkown" typo -> "cep290_unknown" on the genotype sheet ---
$ws = $wb.Worksheets.Item("genotype")

$cellsToFix = @("I4","J5","K5","F6","J6","K6","L6","M6","D9","E9","I9")
foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = "cep290_unknown"
}

# Select the new active cell on the genotype sheet and make it the active tab
$ws.Activate()
$ws.Range("D34").Select()

# --- start_age_hpf sheet loses tab-selection (no cell-selection change) ---
$ws2 = $wb.Worksheets.Item("start_age_hpf")
$ws2.Range("B2:M9").Select()

# Re-activate genotype last so it ends up as the active/selected sheet tab
$ws.Activate()
